$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:G14 (column order: B, C, D, E, F, G)
$data = @{
    2  = @(1.455362044514542, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 1, 25.99591228164478)
    3  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    4  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 8.974608811992548)
    5  = @(0.1190320826869504, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1, 1.406728370586922)
    6  = @(0.6606524410359556, 0.002571899574220771, 0.7527432677738641, 10.19245300693656, 1, 11.6084206153206)
    7  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    8  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    9  = @(3.286832544864788, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1, 4.23731228292506)
    10 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 10.19245300693656, 1, 16.0466451790395)
    11 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    12 = @(0.04271373187048222, 1.655778082260271, 22.3905356188092, 1133.036916526867, 1, 1157.125943959807)
    13 = @(0.2917716402565462, 10.34677158129881, 0.1494219747398047, 1133.036916526867, 1, 1143.824881723162)
    14 = @(0.01293466051926884, 0.306821227259698, 22.3905356188092, 2195978.878461985, 1, 2196001.588753492)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
